# Append 45 new transaction rows (r=1223..1267) to the "Konto" sheet,
# covering the week ending 2021-06-20 (matches the author commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1223
$ws.Cells.Item(1223, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1223, 1).Value = 44361
$ws.Cells.Item(1223, 2).Value = "Reko224"
$ws.Cells.Item(1223, 3).Value = 3011
$ws.Cells.Item(1223, 4).Value = "Reko Swish +46704564448"
$ws.Cells.Item(1223, 6).Value = 345.54

# Row 1224
$ws.Cells.Item(1224, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1224, 1).Value = 44361
$ws.Cells.Item(1224, 2).Value = "Reko224"
$ws.Cells.Item(1224, 3).Value = 2611
$ws.Cells.Item(1224, 4).Value = "Reko Swish +46704564448"
$ws.Cells.Item(1224, 6).Value = 41.46

# Row 1225
$ws.Cells.Item(1225, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1225, 1).Value = 44361
$ws.Cells.Item(1225, 2).Value = "Reko224"
$ws.Cells.Item(1225, 3).Value = 1930
$ws.Cells.Item(1225, 4).Value = "Reko Swish +46704564448"
$ws.Cells.Item(1225, 5).Value = 387

# Row 1226
$ws.Cells.Item(1226, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1226, 1).Value = 44361
$ws.Cells.Item(1226, 2).Value = "Reko225"
$ws.Cells.Item(1226, 3).Value = 3011
$ws.Cells.Item(1226, 4).Value = "Reko Swish +46739806822"
$ws.Cells.Item(1226, 6).Value = 230.36

# Row 1227
$ws.Cells.Item(1227, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1227, 1).Value = 44361
$ws.Cells.Item(1227, 2).Value = "Reko225"
$ws.Cells.Item(1227, 3).Value = 2611
$ws.Cells.Item(1227, 4).Value = "Reko Swish +46739806822"
$ws.Cells.Item(1227, 6).Value = 27.64

# Row 1228
$ws.Cells.Item(1228, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1228, 1).Value = 44361
$ws.Cells.Item(1228, 2).Value = "Reko225"
$ws.Cells.Item(1228, 3).Value = 1930
$ws.Cells.Item(1228, 4).Value = "Reko Swish +46739806822"
$ws.Cells.Item(1228, 5).Value = 258

# Row 1229
$ws.Cells.Item(1229, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1229, 1).Value = 44361
$ws.Cells.Item(1229, 2).Value = "Reko226"
$ws.Cells.Item(1229, 3).Value = 3011
$ws.Cells.Item(1229, 4).Value = "Reko Swish +46709334300"
$ws.Cells.Item(1229, 6).Value = 141.07

# Row 1230
$ws.Cells.Item(1230, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1230, 1).Value = 44361
$ws.Cells.Item(1230, 2).Value = "Reko226"
$ws.Cells.Item(1230, 3).Value = 2611
$ws.Cells.Item(1230, 4).Value = "Reko Swish +46709334300"
$ws.Cells.Item(1230, 6).Value = 16.93

# Row 1231
$ws.Cells.Item(1231, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1231, 1).Value = 44361
$ws.Cells.Item(1231, 2).Value = "Reko226"
$ws.Cells.Item(1231, 3).Value = 1930
$ws.Cells.Item(1231, 4).Value = "Reko Swish +46709334300"
$ws.Cells.Item(1231, 5).Value = 158

# Row 1232
$ws.Cells.Item(1232, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1232, 1).Value = 44361
$ws.Cells.Item(1232, 2).Value = "Reko227"
$ws.Cells.Item(1232, 3).Value = 3011
$ws.Cells.Item(1232, 4).Value = "Reko Swish +46761910051"
$ws.Cells.Item(1232, 6).Value = 141.07

# Row 1233
$ws.Cells.Item(1233, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1233, 1).Value = 44361
$ws.Cells.Item(1233, 2).Value = "Reko227"
$ws.Cells.Item(1233, 3).Value = 2611
$ws.Cells.Item(1233, 4).Value = "Reko Swish +46761910051"
$ws.Cells.Item(1233, 6).Value = 16.93

# Row 1234
$ws.Cells.Item(1234, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1234, 1).Value = 44361
$ws.Cells.Item(1234, 2).Value = "Reko227"
$ws.Cells.Item(1234, 3).Value = 1930
$ws.Cells.Item(1234, 4).Value = "Reko Swish +46761910051"
$ws.Cells.Item(1234, 5).Value = 158

# Row 1235
$ws.Cells.Item(1235, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1235, 1).Value = 44362
$ws.Cells.Item(1235, 2).NumberFormat = "@"
$ws.Cells.Item(1235, 2).Value = "3151843"
$ws.Cells.Item(1235, 3).Value = 3011
$ws.Cells.Item(1235, 4).Value = "Order 3151843 Card(Stripe)"
$ws.Cells.Item(1235, 6).Value = 1912.5

# Row 1236
$ws.Cells.Item(1236, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1236, 1).Value = 44362
$ws.Cells.Item(1236, 2).NumberFormat = "@"
$ws.Cells.Item(1236, 2).Value = "3151843"
$ws.Cells.Item(1236, 3).Value = 2611
$ws.Cells.Item(1236, 4).Value = "Order 3151843 Card(Stripe)"
$ws.Cells.Item(1236, 6).Value = 229.5

# Row 1237
$ws.Cells.Item(1237, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1237, 1).Value = 44362
$ws.Cells.Item(1237, 2).NumberFormat = "@"
$ws.Cells.Item(1237, 2).Value = "3151843"
$ws.Cells.Item(1237, 3).Value = 1930
$ws.Cells.Item(1237, 4).Value = "Order 3151843 Card(Stripe)"
$ws.Cells.Item(1237, 5).Value = 2142

# Row 1238
$ws.Cells.Item(1238, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1238, 1).Value = 44362
$ws.Cells.Item(1238, 2).Value = "Reko228"
$ws.Cells.Item(1238, 3).Value = 3011
$ws.Cells.Item(1238, 4).Value = "Reko Swish +46762568457"
$ws.Cells.Item(1238, 6).Value = 70.54000000000001

# Row 1239
$ws.Cells.Item(1239, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1239, 1).Value = 44362
$ws.Cells.Item(1239, 2).Value = "Reko228"
$ws.Cells.Item(1239, 3).Value = 2611
$ws.Cells.Item(1239, 4).Value = "Reko Swish +46762568457"
$ws.Cells.Item(1239, 6).Value = 8.460000000000001

# Row 1240
$ws.Cells.Item(1240, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1240, 1).Value = 44362
$ws.Cells.Item(1240, 2).Value = "Reko228"
$ws.Cells.Item(1240, 3).Value = 1930
$ws.Cells.Item(1240, 4).Value = "Reko Swish +46762568457"
$ws.Cells.Item(1240, 5).Value = 79

# Row 1241
$ws.Cells.Item(1241, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1241, 1).Value = 44362
$ws.Cells.Item(1241, 3).Value = 4010
$ws.Cells.Item(1241, 4).Value = "TINGSTAD PAPPER"
$ws.Cells.Item(1241, 5).Value = 4156.8

# Row 1242
$ws.Cells.Item(1242, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1242, 1).Value = 44362
$ws.Cells.Item(1242, 3).Value = 2641
$ws.Cells.Item(1242, 4).Value = "TINGSTAD PAPPER"
$ws.Cells.Item(1242, 5).Value = 1039.2

# Row 1243
$ws.Cells.Item(1243, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1243, 1).Value = 44362
$ws.Cells.Item(1243, 3).Value = 1930
$ws.Cells.Item(1243, 4).Value = "TINGSTAD PAPPER"
$ws.Cells.Item(1243, 6).Value = 5196

# Row 1244
$ws.Cells.Item(1244, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1244, 1).Value = 44362
$ws.Cells.Item(1244, 3).Value = 4010
$ws.Cells.Item(1244, 4).Value = "M&S RB BROMMA K0135"
$ws.Cells.Item(1244, 5).Value = 317.2

# Row 1245
$ws.Cells.Item(1245, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1245, 1).Value = 44362
$ws.Cells.Item(1245, 3).Value = 2645
$ws.Cells.Item(1245, 4).Value = "M&S RB BROMMA K0135"
$ws.Cells.Item(1245, 5).Value = 38.06

# Row 1246
$ws.Cells.Item(1246, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1246, 1).Value = 44362
$ws.Cells.Item(1246, 3).Value = 1930
$ws.Cells.Item(1246, 4).Value = "M&S RB BROMMA K0135"
$ws.Cells.Item(1246, 6).Value = 355.26

# Row 1247
$ws.Cells.Item(1247, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1247, 1).Value = 44363
$ws.Cells.Item(1247, 2).Value = "Reko229"
$ws.Cells.Item(1247, 3).Value = 4010
$ws.Cells.Item(1247, 4).Value = "Reko Swish +46704564448 Return"
$ws.Cells.Item(1247, 5).Value = 44.64

# Row 1248
$ws.Cells.Item(1248, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1248, 1).Value = 44363
$ws.Cells.Item(1248, 2).Value = "Reko229"
$ws.Cells.Item(1248, 3).Value = 2645
$ws.Cells.Item(1248, 4).Value = "Reko Swish +46704564448 Return"
$ws.Cells.Item(1248, 5).Value = 5.36

# Row 1249
$ws.Cells.Item(1249, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1249, 1).Value = 44363
$ws.Cells.Item(1249, 2).Value = "Reko229"
$ws.Cells.Item(1249, 3).Value = 1930
$ws.Cells.Item(1249, 4).Value = "Reko Swish +46704564448 Return"
$ws.Cells.Item(1249, 6).Value = 50

# Row 1250
$ws.Cells.Item(1250, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1250, 1).Value = 44363
$ws.Cells.Item(1250, 2).Value = "Reko230"
$ws.Cells.Item(1250, 3).Value = 3011
$ws.Cells.Item(1250, 4).Value = "Reko Swish +46703677212"
$ws.Cells.Item(1250, 6).Value = 460.71

# Row 1251
$ws.Cells.Item(1251, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1251, 1).Value = 44363
$ws.Cells.Item(1251, 2).Value = "Reko230"
$ws.Cells.Item(1251, 3).Value = 2611
$ws.Cells.Item(1251, 4).Value = "Reko Swish +46703677212"
$ws.Cells.Item(1251, 6).Value = 55.29

# Row 1252
$ws.Cells.Item(1252, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1252, 1).Value = 44363
$ws.Cells.Item(1252, 2).Value = "Reko230"
$ws.Cells.Item(1252, 3).Value = 1930
$ws.Cells.Item(1252, 4).Value = "Reko Swish +46703677212"
$ws.Cells.Item(1252, 5).Value = 516

# Row 1253
$ws.Cells.Item(1253, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1253, 1).Value = 44364
$ws.Cells.Item(1253, 2).NumberFormat = "@"
$ws.Cells.Item(1253, 2).Value = "0172055"
$ws.Cells.Item(1253, 3).Value = 3011
$ws.Cells.Item(1253, 4).Value = "Order 0172055 Swish +46703564388"
$ws.Cells.Item(1253, 6).Value = 1062.5

# Row 1254
$ws.Cells.Item(1254, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1254, 1).Value = 44364
$ws.Cells.Item(1254, 2).NumberFormat = "@"
$ws.Cells.Item(1254, 2).Value = "0172055"
$ws.Cells.Item(1254, 3).Value = 2611
$ws.Cells.Item(1254, 4).Value = "Order 0172055 Swish +46703564388"
$ws.Cells.Item(1254, 6).Value = 127.5

# Row 1255
$ws.Cells.Item(1255, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1255, 1).Value = 44364
$ws.Cells.Item(1255, 2).NumberFormat = "@"
$ws.Cells.Item(1255, 2).Value = "0172055"
$ws.Cells.Item(1255, 3).Value = 1930
$ws.Cells.Item(1255, 4).Value = "Order 0172055 Swish +46703564388"
$ws.Cells.Item(1255, 5).Value = 1190

# Row 1256
$ws.Cells.Item(1256, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1256, 1).Value = 44365
$ws.Cells.Item(1256, 2).Value = "Reko231"
$ws.Cells.Item(1256, 3).Value = 3011
$ws.Cells.Item(1256, 4).Value = "Reko Swish +46706183071"
$ws.Cells.Item(1256, 6).Value = 460.71

# Row 1257
$ws.Cells.Item(1257, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1257, 1).Value = 44365
$ws.Cells.Item(1257, 2).Value = "Reko231"
$ws.Cells.Item(1257, 3).Value = 2611
$ws.Cells.Item(1257, 4).Value = "Reko Swish +46706183071"
$ws.Cells.Item(1257, 6).Value = 55.29

# Row 1258
$ws.Cells.Item(1258, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1258, 1).Value = 44365
$ws.Cells.Item(1258, 2).Value = "Reko231"
$ws.Cells.Item(1258, 3).Value = 1930
$ws.Cells.Item(1258, 4).Value = "Reko Swish +46706183071"
$ws.Cells.Item(1258, 5).Value = 516

# Row 1259
$ws.Cells.Item(1259, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1259, 1).Value = 44365
$ws.Cells.Item(1259, 2).Value = "Reko232"
$ws.Cells.Item(1259, 3).Value = 3011
$ws.Cells.Item(1259, 4).Value = "Reko Swish +46722208030"
$ws.Cells.Item(1259, 6).Value = 230.36

# Row 1260
$ws.Cells.Item(1260, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1260, 1).Value = 44365
$ws.Cells.Item(1260, 2).Value = "Reko232"
$ws.Cells.Item(1260, 3).Value = 2611
$ws.Cells.Item(1260, 4).Value = "Reko Swish +46722208030"
$ws.Cells.Item(1260, 6).Value = 27.64

# Row 1261
$ws.Cells.Item(1261, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1261, 1).Value = 44365
$ws.Cells.Item(1261, 2).Value = "Reko232"
$ws.Cells.Item(1261, 3).Value = 1930
$ws.Cells.Item(1261, 4).Value = "Reko Swish +46722208030"
$ws.Cells.Item(1261, 5).Value = 258

# Row 1262
$ws.Cells.Item(1262, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1262, 1).Value = 44365
$ws.Cells.Item(1262, 3).Value = 5670
$ws.Cells.Item(1262, 4).Value = "ST1 V#LLINGBY K0135"
$ws.Cells.Item(1262, 5).Value = 841.16

# Row 1263
$ws.Cells.Item(1263, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1263, 1).Value = 44365
$ws.Cells.Item(1263, 3).Value = 2641
$ws.Cells.Item(1263, 4).Value = "ST1 V#LLINGBY K0135"
$ws.Cells.Item(1263, 5).Value = 210.29

# Row 1264
$ws.Cells.Item(1264, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1264, 1).Value = 44365
$ws.Cells.Item(1264, 3).Value = 1930
$ws.Cells.Item(1264, 4).Value = "ST1 V#LLINGBY K0135"
$ws.Cells.Item(1264, 6).Value = 1051.45

# Row 1265
$ws.Cells.Item(1265, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1265, 1).Value = 44365
$ws.Cells.Item(1265, 3).Value = 4010
$ws.Cells.Item(1265, 4).Value = "M&S RB BROMMA K0135"
$ws.Cells.Item(1265, 5).Value = 529.4

# Row 1266
$ws.Cells.Item(1266, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1266, 1).Value = 44365
$ws.Cells.Item(1266, 3).Value = 2645
$ws.Cells.Item(1266, 4).Value = "M&S RB BROMMA K0135"
$ws.Cells.Item(1266, 5).Value = 63.53

# Row 1267
$ws.Cells.Item(1267, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1267, 1).Value = 44365
$ws.Cells.Item(1267, 3).Value = 1930
$ws.Cells.Item(1267, 4).Value = "M&S RB BROMMA K0135"
$ws.Cells.Item(1267, 6).Value = 592.9299999999999

